$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 45415.043
$ws.Range("I19").Value = 100737.1
$ws.Range("J19").Value = 2859.6155
$ws.Range("K19").Value = 100737.1
$ws.Range("L19").Value = 2859.6155
$ws.Range("M19").Value = -100562.1
$ws.Range("N19").Value = -3209.6155
$ws.Range("H40").Value = 7248.75
$ws.Range("I40").Value = 3777.4
$ws.Range("K40").Value = 3777.4
$ws.Range("M40").Value = -3602.4
$ws.Range("H82").Value = 11373.286
$ws.Range("I82").Value = 11373.286
$ws.Range("K82").Value = 34119.858
$ws.Range("M82").Value = -33713.858
$ws.Range("H85").Value = 11373.286
$ws.Range("I85").Value = 11373.286
$ws.Range("K85").Value = 34119.858
$ws.Range("M85").Value = -32715.858
$ws.Range("H125").Value = 98852.73
$ws.Range("I125").Value = 992.5
$ws.Range("J125").Value = 120599.445
$ws.Range("K125").Value = 8932.5
$ws.Range("L125").Value = 1085395.005
$ws.Range("M125").Value = -6472.5
$ws.Range("N125").Value = -1090315.005
$ws.Range("H131").Value = 6906.7144
$ws.Range("I131").Value = 1800
$ws.Range("K131").Value = 5400
$ws.Range("M131").Value = -360
$ws.Range("H141").Value = 5340.8
$ws.Range("I141").Value = 2234.6667
$ws.Range("K141").Value = 6704.000100000001
$ws.Range("M141").Value = -1524.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2833.8591
$ws.Range("I32").Value = 2231.9092
$ws.Range("J32").Value = 10779.6
$ws.Range("K32").Value = 2231.9092
$ws.Range("L32").Value = 10779.6
$ws.Range("M32").Value = -1944.9092
$ws.Range("N32").Value = -11353.6
$ws.Range("H45").Value = 2203.5334
$ws.Range("I45").Value = 1421.75
$ws.Range("K45").Value = 1421.75
$ws.Range("M45").Value = -1044.75
$ws.Range("H102").Value = 2183.7827
$ws.Range("I102").Value = 1838.85
$ws.Range("J102").Value = 4483.3335
$ws.Range("K102").Value = 1838.85
$ws.Range("L102").Value = 4483.3335
$ws.Range("M102").Value = -216.8499999999999
$ws.Range("N102").Value = -7727.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1976.1923
$ws.Range("I107").Value = 2097.238
$ws.Range("K107").Value = 2097.238
$ws.Range("M107").Value = -177.2379999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2199.75
$ws.Range("I16").Value = 1749.5
$ws.Range("K16").Value = 1749.5
$ws.Range("M16").Value = -1462.5
$ws.Range("H31").Value = 6870.9395
$ws.Range("I31").Value = 2992.4736
$ws.Range("J31").Value = 12134.571
$ws.Range("K31").Value = 2992.4736
$ws.Range("L31").Value = 12134.571
$ws.Range("M31").Value = -2697.4736
$ws.Range("N31").Value = -12724.571
$ws.Range("H34").Value = 6870.9395
$ws.Range("I34").Value = 2992.4736
$ws.Range("J34").Value = 12134.571
$ws.Range("K34").Value = 2992.4736
$ws.Range("L34").Value = 12134.571
$ws.Range("M34").Value = -2790.4736
$ws.Range("N34").Value = -12538.571
$ws.Range("H58").Value = 5758.5884
$ws.Range("I58").Value = 2999.818
$ws.Range("K58").Value = 2999.818
$ws.Range("M58").Value = -2796.818
$ws.Range("H62").Value = 4673.1177
$ws.Range("I62").Value = 4418.615
$ws.Range("J62").Value = 5500.25
$ws.Range("K62").Value = 4418.615
$ws.Range("L62").Value = 5500.25
$ws.Range("M62").Value = -3794.615
$ws.Range("N62").Value = -6748.25
$ws.Range("H65").Value = 4673.1177
$ws.Range("I65").Value = 4418.615
$ws.Range("J65").Value = 5500.25
$ws.Range("K65").Value = 22093.075
$ws.Range("L65").Value = 27501.25
$ws.Range("M65").Value = -18973.075
$ws.Range("N65").Value = -33741.25
$ws.Range("H94").Value = 2476.6
$ws.Range("I94").Value = 1477.25
$ws.Range("J94").Value = 3618.7144
$ws.Range("K94").Value = 1477.25
$ws.Range("L94").Value = 3618.7144
$ws.Range("M94").Value = -1026.25
$ws.Range("N94").Value = -4520.7144
$ws.Range("H105").Value = 2313.5557
$ws.Range("I105").Value = 1979.55
$ws.Range("K105").Value = 1979.55
$ws.Range("M105").Value = -232.55
$ws.Range("H107").Value = 1137
$ws.Range("I107").Value = 692.2857
$ws.Range("K107").Value = 692.2857
$ws.Range("M107").Value = 1227.7143
$ws.Range("H113").Value = 2199.75
$ws.Range("I113").Value = 1749.5
$ws.Range("K113").Value = 1749.5
$ws.Range("M113").Value = 420.5
$ws.Range("H132").Value = 2752.2766
$ws.Range("I132").Value = 2241.5476
$ws.Range("J132").Value = 7042.4
$ws.Range("K132").Value = 6724.6428
$ws.Range("L132").Value = 21127.2
$ws.Range("M132").Value = -4194.6428
$ws.Range("N132").Value = -26187.2
$ws.Range("H134").Value = 6458.8696
$ws.Range("I134").Value = 3644.4443
$ws.Range("K134").Value = 10933.3329
$ws.Range("M134").Value = -8398.332900000001
$ws.Range("H136").Value = 5758.5884
$ws.Range("I136").Value = 2999.818
$ws.Range("K136").Value = 8999.454000000002
$ws.Range("M136").Value = -6449.454000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 80.09524
$ws.Range("I2").Value = 102.92857
$ws.Range("J2").Value = 34.42857
$ws.Range("K2").Value = 617.57142
$ws.Range("L2").Value = 206.57142
$ws.Range("M2").Value = -504.57142
$ws.Range("N2").Value = -432.57142
$ws.Range("H15").Value = 221.66667
$ws.Range("I15").Value = 166
$ws.Range("K15").Value = 498
$ws.Range("M15").Value = -358
$ws.Range("H17").Value = 341
$ws.Range("I17").Value = 314.33334
$ws.Range("J17").Value = 357
$ws.Range("K17").Value = 943.0000200000001
$ws.Range("L17").Value = 1071
$ws.Range("M17").Value = -774.0000200000001
$ws.Range("N17").Value = -1409
$ws.Range("H33").Value = 330.63635
$ws.Range("I33").Value = 219.06667
$ws.Range("K33").Value = 1314.40002
$ws.Range("M33").Value = -1031.40002
$ws.Range("H34").Value = 3422.7273
$ws.Range("I34").Value = 225
$ws.Range("J34").Value = 4133.3335
$ws.Range("K34").Value = 675
$ws.Range("L34").Value = 12400.0005
$ws.Range("M34").Value = -591
$ws.Range("N34").Value = -12568.0005
$ws.Range("H38").Value = 41.46154
$ws.Range("I38").Value = 42.25
$ws.Range("K38").Value = 126.75
$ws.Range("M38").Value = 220.25
$ws.Range("H97").Value = 4063.7646
$ws.Range("I97").Value = 3156.3333
$ws.Range("J97").Value = 4558.727
$ws.Range("K97").Value = 9468.999899999999
$ws.Range("L97").Value = 13676.181
$ws.Range("M97").Value = -8972.999899999999
$ws.Range("N97").Value = -14668.181
$ws.Range("H98").Value = 3181.0908
$ws.Range("I98").Value = 2241.1667
$ws.Range("J98").Value = 4309
$ws.Range("K98").Value = 6723.500100000001
$ws.Range("L98").Value = 12927
$ws.Range("M98").Value = -5225.500100000001
$ws.Range("N98").Value = -15923

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1154.5217
$ws.Range("I2").Value = 788.0952
$ws.Range("J2").Value = 5002
$ws.Range("K2").Value = 788.0952
$ws.Range("L2").Value = 5002
$ws.Range("M2").Value = -675.0952
$ws.Range("N2").Value = -5228
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H80").Value = 4814.684
$ws.Range("I80").Value = 4168.7144
$ws.Range("J80").Value = 5191.5
$ws.Range("K80").Value = 4168.7144
$ws.Range("L80").Value = 5191.5
$ws.Range("M80").Value = -3170.7144
$ws.Range("N80").Value = -7187.5
$ws.Range("H83").Value = 4814.684
$ws.Range("I83").Value = 4168.7144
$ws.Range("J83").Value = 5191.5
$ws.Range("K83").Value = 20843.572
$ws.Range("L83").Value = 25957.5
$ws.Range("M83").Value = -15851.572
$ws.Range("N83").Value = -35941.5
$ws.Range("H132").Value = 3979
$ws.Range("I132").Value = 3075.7368
$ws.Range("K132").Value = 9227.2104
$ws.Range("M132").Value = -6697.2104

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1527.36
$ws.Range("I16").Value = 787.6087
$ws.Range("K16").Value = 787.6087
$ws.Range("M16").Value = -617.6087
$ws.Range("H61").Value = 3835.2896
$ws.Range("I61").Value = 3080.9
$ws.Range("K61").Value = 3080.9
$ws.Range("M61").Value = -2878.9
$ws.Range("H109").Value = 175884.17
$ws.Range("J109").Value = 175884.17
$ws.Range("L109").Value = 175884.17
$ws.Range("N109").Value = -178658.17
$ws.Range("H113").Value = 3835.2896
$ws.Range("I113").Value = 3080.9
$ws.Range("K113").Value = 3080.9
$ws.Range("M113").Value = -910.9000000000001
$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -59840
$ws.Range("H136").Value = 8031.2666
$ws.Range("I136").Value = 4400.3335
$ws.Range("K136").Value = 13201.0005
$ws.Range("M136").Value = -10651.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1347.4333
$ws.Range("I107").Value = 1268.1364
$ws.Range("J107").Value = 1565.5
$ws.Range("K107").Value = 3804.4092
$ws.Range("L107").Value = 4696.5
$ws.Range("M107").Value = -1884.4092
$ws.Range("N107").Value = -8536.5
$ws.Range("H126").Value = 3982.8333
$ws.Range("I126").Value = 2974.25
$ws.Range("K126").Value = 8922.75
$ws.Range("M126").Value = -6452.75
$ws.Range("H132").Value = 5628.8335
$ws.Range("I132").Value = 4104.909
$ws.Range("K132").Value = 12314.727
$ws.Range("M132").Value = -9784.726999999999
$ws.Range("H136").Value = 4472.05
$ws.Range("I136").Value = 3876
$ws.Range("J136").Value = 8644.4
$ws.Range("K136").Value = 11628
$ws.Range("L136").Value = 25933.2
$ws.Range("M136").Value = -9078
$ws.Range("N136").Value = -31033.2
